$d = $word.ActiveDocument

# 1. In the "Time Series Plot" paragraph, drop the trailing space after
#    "...Y axis." and append the new sentence about the time series plots.
$null = $d.Content.Find.Execute( `
    "plotted on the Y axis. ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "plotted on the Y axis. Below are the time series plot of the input data and the output data. ", `
    2)

# 2. Add a new (mostly blank) "Paragraph"-styled paragraph right after it,
#    containing just two spaces.
$target = $d.Paragraphs.Last
$target.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$newPara.Range.InsertBefore("  ")
